$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: shift existing data block (rows 184-282, columns A-R) down by 4 rows
# to make room for 4 new rows of data (new rows 184-187).
$src = $ws.Range("A184:R282")
$src.Copy()
$dest = $ws.Range("A188")
$dest.PasteSpecial()
$excel.CutCopyMode = $false

# Step 2: populate the 4 newly freed rows (184-187) with the new records.

# Row 184
$ws.Cells.Item(184, 1).Value = 10
$ws.Cells.Item(184, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(184, 3).Value = "La Araucanía"
$ws.Cells.Item(184, 4).Value = 44452
$ws.Cells.Item(184, 5).Value = 9
$ws.Cells.Item(184, 6).Value = 100112043
$ws.Cells.Item(184, 7).Value = "Pepino ensalada"
$ws.Cells.Item(184, 8).Value = "Alaska"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 80
$ws.Cells.Item(184, 11).Value = 20000
$ws.Cells.Item(184, 12).Value = 20000
$ws.Cells.Item(184, 13).Value = 20000
$ws.Cells.Item(184, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(184, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(184, 16).Value = 333
$ws.Cells.Item(184, 17).Value = 60
$ws.Cells.Item(184, 18).Value = "Hortaliza"

# Row 185
$ws.Cells.Item(185, 1).Value = 10
$ws.Cells.Item(185, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(185, 3).Value = "La Araucanía"
$ws.Cells.Item(185, 4).Value = 44452
$ws.Cells.Item(185, 5).Value = 9
$ws.Cells.Item(185, 6).Value = 100112043
$ws.Cells.Item(185, 7).Value = "Pepino ensalada"
$ws.Cells.Item(185, 8).Value = "Sin especificar"
$ws.Cells.Item(185, 9).Value = "Extra"
$ws.Cells.Item(185, 10).Value = 90
$ws.Cells.Item(185, 11).Value = 19000
$ws.Cells.Item(185, 12).Value = 20000
$ws.Cells.Item(185, 13).Value = 19556
$ws.Cells.Item(185, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(185, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(185, 16).Value = 391
$ws.Cells.Item(185, 17).Value = 50
$ws.Cells.Item(185, 18).Value = "Hortaliza"

# Row 186
$ws.Cells.Item(186, 1).Value = 10
$ws.Cells.Item(186, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(186, 3).Value = "La Araucanía"
$ws.Cells.Item(186, 4).Value = 44452
$ws.Cells.Item(186, 5).Value = 9
$ws.Cells.Item(186, 6).Value = 100112043
$ws.Cells.Item(186, 7).Value = "Pepino ensalada"
$ws.Cells.Item(186, 8).Value = "Sin especificar"
$ws.Cells.Item(186, 9).Value = "Primera"
$ws.Cells.Item(186, 10).Value = 600
$ws.Cells.Item(186, 11).Value = 16000
$ws.Cells.Item(186, 12).Value = 17000
$ws.Cells.Item(186, 13).Value = 16500
$ws.Cells.Item(186, 14).Value = "`$/caja 60 unidades"
$ws.Cells.Item(186, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(186, 16).Value = 275
$ws.Cells.Item(186, 17).Value = 60
$ws.Cells.Item(186, 18).Value = "Hortaliza"

# Row 187
$ws.Cells.Item(187, 1).Value = 10
$ws.Cells.Item(187, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(187, 3).Value = "La Araucanía"
$ws.Cells.Item(187, 4).Value = 44452
$ws.Cells.Item(187, 5).Value = 9
$ws.Cells.Item(187, 6).Value = 100112043
$ws.Cells.Item(187, 7).Value = "Pepino ensalada"
$ws.Cells.Item(187, 8).Value = "Sin especificar"
$ws.Cells.Item(187, 9).Value = "Segunda"
$ws.Cells.Item(187, 10).Value = 90
$ws.Cells.Item(187, 11).Value = 12000
$ws.Cells.Item(187, 12).Value = 14000
$ws.Cells.Item(187, 13).Value = 13111
$ws.Cells.Item(187, 14).Value = "`$/caja 80 unidades"
$ws.Cells.Item(187, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(187, 16).Value = 164
$ws.Cells.Item(187, 17).Value = 80
$ws.Cells.Item(187, 18).Value = "Hortaliza"

# Update the sheet dimension reference to reflect the new extent.
$ws.Range("A1:R286").Select()
